$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17: One for the Road / Potion
$ws.Range("H17").Value = 60273.06
$ws.Range("J17").Value = 60273.06
$ws.Range("L17").Value = 180819.18
$ws.Range("N17").Value = -181155.18

# Row 41: The Write Stuff / Enchanted Mythril Ink
$ws.Range("H41").Value = 3966.3333
$ws.Range("I41").Value = 4819.5
$ws.Range("K41").Value = 4819.5
$ws.Range("M41").Value = -4379.5

# Row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 9399.143
$ws.Range("I76").Value = 10299
$ws.Range("K76").Value = 10299
$ws.Range("M76").Value = -9984

# Row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 9399.143
$ws.Range("I79").Value = 10299
$ws.Range("K79").Value = 10299
$ws.Range("M79").Value = -9207

# Row 106: Making Your Mark / Enchanted Palladium Ink
$ws.Range("H106").Value = 16849.834
$ws.Range("I106").Value = 16849.834
$ws.Range("K106").Value = 16849.834
$ws.Range("M106").Value = -16218.834

# Row 113: Amaro Kart / Starch Glue
$ws.Range("H113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = $null

# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 5734.1313
$ws.Range("I138").Value = 4781
$ws.Range("J138").Value = 5992.271
$ws.Range("K138").Value = 14343
$ws.Range("L138").Value = 17976.813
$ws.Range("M138").Value = -9203
$ws.Range("N138").Value = -28256.813

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 1892.3334
$ws.Range("I2").Value = 2088.5
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 2088.5
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -1975.5
$ws.Range("N2").Value = -1726

# Row 46: Get Me the Usual / Heavy Steel Flanchard
$ws.Range("H46").Value = 19830.8
$ws.Range("I46").Value = 10979.333
$ws.Range("K46").Value = 10979.333
$ws.Range("M46").Value = -10660.333

# Row 63: Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 3548.3333
$ws.Range("I63").Value = 3548.3333
$ws.Range("K63").Value = 3548.3333
$ws.Range("M63").Value = -2862.3333

# Row 66: A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 3548.3333
$ws.Range("I66").Value = 3548.3333
$ws.Range("K66").Value = 17741.6665
$ws.Range("M66").Value = -14309.6665

# Row 88: The Mast Chance / Adamantite Rivets
$ws.Range("H88").Value = 3733.1667
$ws.Range("I88").Value = 2700
$ws.Range("K88").Value = 2700
$ws.Range("M88").Value = -2294

# Row 91: The Rose and the Riveter (L) / Adamantite Rivets
$ws.Range("H91").Value = 3733.1667
$ws.Range("I91").Value = 2700
$ws.Range("K91").Value = 2700
$ws.Range("M91").Value = -1296

# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 1892.3334
$ws.Range("I116").Value = 2088.5
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 2088.5
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = 205.5
$ws.Range("N116").Value = -6088

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 2967.919
$ws.Range("I132").Value = 3179.4827
$ws.Range("J132").Value = 2201
$ws.Range("K132").Value = 9538.4481
$ws.Range("L132").Value = 6603
$ws.Range("M132").Value = -7008.4481
$ws.Range("N132").Value = -11663

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 1892.3334
$ws.Range("I3").Value = 2088.5
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 2088.5
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = -1974.5
$ws.Range("N3").Value = -1728

# Row 64: With Bearings Straight / Mythrite Nugget
$ws.Range("H64").Value = 19231346
$ws.Range("I64").Value = 62500452
$ws.Range("J64").Value = 632
$ws.Range("K64").Value = 62500452
$ws.Range("L64").Value = 632
$ws.Range("M64").Value = -62500227
$ws.Range("N64").Value = -1082

# Row 67: Bearing the Brunt (L) / Mythrite Nugget
$ws.Range("H67").Value = 19231346
$ws.Range("I67").Value = 62500452
$ws.Range("J67").Value = 632
$ws.Range("K67").Value = 62500452
$ws.Range("L67").Value = 632
$ws.Range("M67").Value = -62499672
$ws.Range("N67").Value = -2192

# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 647470.0600000001
$ws.Range("I105").Value = 859839.5600000001
$ws.Range("K105").Value = 859839.5600000001
$ws.Range("M105").Value = -858092.5600000001

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 3417.5
$ws.Range("I134").Value = 2339
$ws.Range("K134").Value = 7017
$ws.Range("M134").Value = -4482

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 20003124
$ws.Range("I31").Value = 25643682
$ws.Range("K31").Value = 25643682
$ws.Range("M31").Value = -25643387

# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 20003124
$ws.Range("I34").Value = 25643682
$ws.Range("K34").Value = 25643682
$ws.Range("M34").Value = -25643480

# Row 43: The Long Lance of the Law / Steel Halberd
$ws.Range("H43").Value = 18216
$ws.Range("J43").Value = 18216
$ws.Range("L43").Value = 18216
$ws.Range("N43").Value = -18584

# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 21285.637
$ws.Range("I62").Value = 14749
$ws.Range("J62").Value = 26732.834
$ws.Range("K62").Value = 14749
$ws.Range("L62").Value = 26732.834
$ws.Range("M62").Value = -14125
$ws.Range("N62").Value = -27980.834

# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 21285.637
$ws.Range("I65").Value = 14749
$ws.Range("J65").Value = 26732.834
$ws.Range("K65").Value = 73745
$ws.Range("L65").Value = 133664.17
$ws.Range("M65").Value = -70625
$ws.Range("N65").Value = -139904.17

# Row 101: Everybody's Heard about the 'Berd / Doman Steel Halberd
$ws.Range("H101").Value = 18216
$ws.Range("J101").Value = 18216
$ws.Range("L101").Value = 18216
$ws.Range("N101").Value = -24706

$ws = $wb.Worksheets.Item("CUL")
# Row 121: A Cookie for Your Troubles / Coffee Biscuit
$ws.Range("H121").Value = 2653.2942
$ws.Range("I121").Value = 453.8
$ws.Range("K121").Value = 1361.4
$ws.Range("M121").Value = -51.40000000000009

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 5185.4287
$ws.Range("I70").Value = 3498
$ws.Range("K70").Value = 3498
$ws.Range("M70").Value = -3228

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 5185.4287
$ws.Range("I73").Value = 3498
$ws.Range("K73").Value = 3498
$ws.Range("M73").Value = -2562

# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 4481.5625
$ws.Range("I80").Value = 3669.3
$ws.Range("J80").Value = 5835.3335
$ws.Range("K80").Value = 3669.3
$ws.Range("L80").Value = 5835.3335
$ws.Range("M80").Value = -2671.3
$ws.Range("N80").Value = -7831.3335

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 4481.5625
$ws.Range("I83").Value = 3669.3
$ws.Range("J83").Value = 5835.3335
$ws.Range("K83").Value = 18346.5
$ws.Range("L83").Value = 29176.6675
$ws.Range("M83").Value = -13354.5
$ws.Range("N83").Value = -39160.6675

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 5000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = $null
$ws.Range("N102").Value = -8244

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 2766.697
$ws.Range("I132").Value = 3095.75
$ws.Range("K132").Value = 9287.25
$ws.Range("M132").Value = -6757.25

$ws = $wb.Worksheets.Item("LTW")
# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 6559.56
$ws.Range("I40").Value = 6064.5713
$ws.Range("K40").Value = 6064.5713
$ws.Range("M40").Value = -5928.5713

$ws = $wb.Worksheets.Item("WVR")
# Row 94: Proper Props / Bloodhempen Armguards of Scouting
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = $null
